# "Modified DDT test login and added cart page object"
# The login DDT sheet gains 4 new data-driven test rows (inserted right
# after the existing "standard_user/ducks/fail" row), pushing the former
# rows 4-7 down to rows 8-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 4 - shifts old rows 4:7 down to 8:11.
$ws.Rows("4:7").Insert() | Out-Null

# Give the new rows the same bordered formatting as the other data rows
# (row 2) by copy/pasting formats only.
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A4:C7").PasteSpecial(-4122) | Out-Null

# New row 4: username / password / fail
$ws.Range("A4").Value = "username"
$ws.Range("B4").Value = "password"
$ws.Range("C4").Value = "fail"

# New row 5: (blank username) / 1234 / fail
$ws.Range("B5").Value = 1234
$ws.Range("C5").Value = "fail"

# New row 6: standard_user / (blank password) / fail
$ws.Range("A6").Value = "standard_user"
$ws.Range("C6").Value = "fail"

# New row 7: (blank username) / (blank password) / fail
$ws.Range("C7").Value = "fail"

# Match the saved selection state from the authored workbook.
$ws.Range("E7").Select() | Out-Null

Write-Output "done"
